# Refactorization: itialize no longer hardcoded
#
# The Clue layout grid on Sheet1 used the literal text "w" (lowercase) as the
# filler/placeholder value baked into most cells. Replace every cell whose
# value is the literal "w" with "W" (uppercase) so the sheet no longer
# depends on that hardcoded lowercase literal.
#
# (Once every "w" cell is repointed to "W", nothing references the old "w"
# shared-string entry any more, so it naturally drops out of the shared
# string table on save, while the new "W" entry is appended.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $used.Cells.Item($r, $c)
        if ($cell.Value2 -ceq "w") {
            $cell.Value = "W"
        }
    }
}

# Leave the selection where the editor last clicked while reviewing the change.
[void]$ws.Range("N23").Select()
